# Apply weekly Fruta/Hortaliza price-sheet update: rows 3-19 (excluding fixed
# row 12) get their varying columns (Fecha, Calidad, Volumen, Precios,
# Unidad, Origen, Precio $/Kg, Kg/unidad) reassigned per the new weekly
# snapshot. Identifying columns (A,B,C,E-K) are untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 <- old row 9
$ws.Range("D3").Value = 44860
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 200
$ws.Range("N3").Value = 23000
$ws.Range("O3").Value = 24000
$ws.Range("P3").Value = 23500
$ws.Range("R3").Value = "Provincia de Limarí"
$ws.Range("S3").Value = 1958

# Row 4 <- old row 14
$ws.Range("D4").Value = 44167
$ws.Range("L4").Value = "Segunda"
$ws.Range("M4").Value = 200
$ws.Range("N4").Value = 18000
$ws.Range("O4").Value = 19000
$ws.Range("P4").Value = 18500
$ws.Range("Q4").Value = "`$/caja 13 kilos"
$ws.Range("S4").Value = 1423
$ws.Range("T4").Value = 13

# Row 5 <- old row 6
$ws.Range("D5").Value = 44545
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 200
$ws.Range("N5").Value = 23000
$ws.Range("O5").Value = 24000
$ws.Range("P5").Value = 23500
$ws.Range("Q5").Value = "`$/bandeja 12 kilos"
$ws.Range("S5").Value = 1958

# Row 6 <- old row 10
$ws.Range("D6").Value = 44783
$ws.Range("L6").Value = "Tercera"
$ws.Range("M6").Value = 100
$ws.Range("N6").Value = 27000
$ws.Range("O6").Value = 28000
$ws.Range("P6").Value = 27500
$ws.Range("Q6").Value = "`$/caja 12 kilos"
$ws.Range("S6").Value = 2292

# Row 7 <- old row 11
$ws.Range("D7").Value = 44496
$ws.Range("N7").Value = 23000
$ws.Range("O7").Value = 24000
$ws.Range("P7").Value = 23500
$ws.Range("Q7").Value = "`$/caja 12 kilos"
$ws.Range("S7").Value = 1958
$ws.Range("T7").Value = 12

# Row 8 <- old row 15
$ws.Range("D8").Value = 44524
$ws.Range("M8").Value = 200
$ws.Range("N8").Value = 23000
$ws.Range("O8").Value = 24000
$ws.Range("P8").Value = 23500
$ws.Range("S8").Value = 1958

# Row 9 <- old row 8
$ws.Range("D9").Value = 44482
$ws.Range("M9").Value = 160
$ws.Range("N9").Value = 25000
$ws.Range("O9").Value = 26000
$ws.Range("P9").Value = 25500
$ws.Range("R9").Value = "Región de Coquimbo"
$ws.Range("S9").Value = 2125

# Row 10 <- old row 19
$ws.Range("D10").Value = 44776
$ws.Range("L10").Value = "Segunda"
$ws.Range("M10").Value = 160
$ws.Range("N10").Value = 29000
$ws.Range("O10").Value = 30000
$ws.Range("P10").Value = 29500
$ws.Range("Q10").Value = "`$/caja 10 kilos"
$ws.Range("S10").Value = 2950
$ws.Range("T10").Value = 10

# Row 11 <- old row 17
$ws.Range("D11").Value = 44811
$ws.Range("M11").Value = 100
$ws.Range("N11").Value = 29000
$ws.Range("O11").Value = 30000
$ws.Range("P11").Value = 29500
$ws.Range("S11").Value = 2458

# Row 13 <- old row 7
$ws.Range("D13").Value = 44468
$ws.Range("L13").Value = "Primera"
$ws.Range("M13").Value = 200
$ws.Range("N13").Value = 29000
$ws.Range("O13").Value = 30000
$ws.Range("P13").Value = 29500
$ws.Range("Q13").Value = "`$/bandeja 10 kilos"
$ws.Range("S13").Value = 2950
$ws.Range("T13").Value = 10

# Row 14 <- old row 3
$ws.Range("D14").Value = 44839
$ws.Range("M14").Value = 160
$ws.Range("N14").Value = 26000
$ws.Range("O14").Value = 27000
$ws.Range("P14").Value = 26500
$ws.Range("Q14").Value = "`$/caja 12 kilos"
$ws.Range("S14").Value = 2208
$ws.Range("T14").Value = 12

# Row 15 <- old row 16
$ws.Range("D15").Value = 44160
$ws.Range("L15").Value = "Segunda"
$ws.Range("N15").Value = 19000
$ws.Range("O15").Value = 20000
$ws.Range("P15").Value = 19500
$ws.Range("Q15").Value = "`$/caja 13 kilos"
$ws.Range("S15").Value = 1500
$ws.Range("T15").Value = 13

# Row 16 <- old row 18
$ws.Range("D16").Value = 44881
$ws.Range("L16").Value = "Primera"
$ws.Range("N16").Value = 22000
$ws.Range("O16").Value = 23000
$ws.Range("P16").Value = 22500
$ws.Range("Q16").Value = "`$/caja 12 kilos"
$ws.Range("S16").Value = 1875
$ws.Range("T16").Value = 12

# Row 17 <- old row 13
$ws.Range("D17").Value = 44874
$ws.Range("L17").Value = "Segunda"
$ws.Range("M17").Value = 250
$ws.Range("N17").Value = 22000
$ws.Range("O17").Value = 23000
$ws.Range("P17").Value = 22500
$ws.Range("S17").Value = 1875

# Row 18 <- old row 4
$ws.Range("D18").Value = 44846
$ws.Range("M18").Value = 160
$ws.Range("N18").Value = 24000
$ws.Range("O18").Value = 25000
$ws.Range("P18").Value = 24500
$ws.Range("S18").Value = 2042

# Row 19 <- old row 5
$ws.Range("D19").Value = 44846
$ws.Range("M19").Value = 100
$ws.Range("N19").Value = 22000
$ws.Range("O19").Value = 23000
$ws.Range("P19").Value = 22500
$ws.Range("Q19").Value = "`$/caja 12 kilos"
$ws.Range("S19").Value = 1875
$ws.Range("T19").Value = 12
